$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: append a new run " (ARMIJO)" right after "...slides 30-35"
# (same paragraph, TASK 2 / main_OPTCON_optimization line)
# ------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(", slides 30-35", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(0)
$text1 = " (ARMIJO)"
$start1 = $rng1.Start
$rng1.InsertAfter($text1)
$rng1.SetRange($start1, $start1 + $text1.Length)
$rng1.Bold = 1
$rng1.LanguageID = "en-GB"
$rng1.Bold = 0

# ------------------------------------------------------------------
# Change 2: append a new run ", slide 36" right after
#           "OPTCON_mpc_2024_11_25, slides 2-3, slide 18" (TASK 4)
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("OPTCON_mpc_2024_11_25, slides 2-3, slide 18", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$text2 = ", slide 36"
$start2 = $rng2.Start
$rng2.InsertAfter($text2)
$rng2.SetRange($start2, $start2 + $text2.Length)
$rng2.Bold = 1
$rng2.LanguageID = "en-GB"
$rng2.Bold = 0

# ------------------------------------------------------------------
# Change 3: brand-new paragraph right after that one with
#           "OPTCON_optimal_constrol_design_2024_11_25, slide 8"
# ------------------------------------------------------------------
$rng2.Collapse(0)
$rng2.InsertParagraphAfter() | Out-Null
$rng2.Collapse(0)
$rng2.MoveStart(1, 1) | Out-Null
$text3 = "OPTCON_optimal_constrol_design_2024_11_25, slide 8"
$rng2.InsertAfter($text3)

Write-Output "edit applied"
